$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.7048620983793006
$ws.Range("C2").Value = 0.6905292479108636
$ws.Range("D2").Value = 0.6976220627550302
$ws.Range("E2").Value = 3590

$ws.Range("B3").Value = 0.4052462526766595
$ws.Range("C3").Value = 0.4217270194986072
$ws.Range("D3").Value = 0.4133224133224133
$ws.Range("E3").Value = 1795

$ws.Range("B4").Value = 0.6009285051067781
$ws.Range("C4").Value = 0.6009285051067781
$ws.Range("D4").Value = 0.6009285051067781
$ws.Range("E4").Value = 0.6009285051067781

$ws.Range("B5").Value = 0.5550541755279801
$ws.Range("C5").Value = 0.5561281337047355
$ws.Range("D5").Value = 0.5554722380387218
$ws.Range("E5").Value = 5385

$ws.Range("B6").Value = 0.6049901498117536
$ws.Range("C6").Value = 0.6009285051067781
$ws.Range("D6").Value = 0.602855512944158
$ws.Range("E6").Value = 5385
